# Commit: "updated test data sheet"
# The SmokeTest column (D) for every scenario row is being switched to "Yes"
# (it was already "Yes" for rows 10-17, and "No" everywhere else).
# The RegressionTest column (E) is left as-is ("No" throughout).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:D335").Value = "Yes"

# Reflect the cursor/selection position left behind by the edit session.
$ws.Range("C341").Select()
